$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting everything down
$ws.Rows.Item(1).Insert()

# Add header values in the new first row
$ws.Range("A1").Value = "Класс"
$ws.Range("B1").Value = "Макс кол-во уроков"

# Update last-selected cell (mirrors the author's working selection)
$ws.Range("F12").Select()
